$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1774.3116
$ws.Range("I15").Value = 1774.3116
$ws.Range("K15").Value = 5322.9348
$ws.Range("M15").Value = -5153.9348
$ws.Range("H17").Value = 33016
$ws.Range("J17").Value = 33016
$ws.Range("L17").Value = 99048
$ws.Range("N17").Value = -99384
$ws.Range("H43").Value = 1301.5
$ws.Range("J43").Value = 1482.4
$ws.Range("L43").Value = 1482.4
$ws.Range("N43").Value = -1620.4
$ws.Range("H74").Value = 4530.6
$ws.Range("I74").Value = 4061.2
$ws.Range("K74").Value = 4061.2
$ws.Range("M74").Value = -3125.2
$ws.Range("H76").Value = 4799.3
$ws.Range("I76").Value = 5300.5
$ws.Range("J76").Value = 4047.5
$ws.Range("K76").Value = 5300.5
$ws.Range("L76").Value = 4047.5
$ws.Range("M76").Value = -4985.5
$ws.Range("N76").Value = -4677.5
$ws.Range("H77").Value = 4530.6
$ws.Range("I77").Value = 4061.2
$ws.Range("K77").Value = 20306
$ws.Range("M77").Value = -15626
$ws.Range("H79").Value = 4799.3
$ws.Range("I79").Value = 5300.5
$ws.Range("J79").Value = 4047.5
$ws.Range("K79").Value = 5300.5
$ws.Range("L79").Value = 4047.5
$ws.Range("M79").Value = -4208.5
$ws.Range("N79").Value = -6231.5
$ws.Range("H96").Value = 723.25
$ws.Range("I96").Value = 456.66666
$ws.Range("J96").Value = 883.2
$ws.Range("K96").Value = 1369.99998
$ws.Range("L96").Value = 2649.6
$ws.Range("M96").Value = 3.00001999999995
$ws.Range("N96").Value = -5395.6
$ws.Range("H121").Value = 1549.0834
$ws.Range("J121").Value = 1713.2222
$ws.Range("L121").Value = 5139.6666
$ws.Range("N121").Value = -8633.6666
$ws.Range("H125").Value = 167149.83
$ws.Range("J125").Value = 400
$ws.Range("L125").Value = 3600
$ws.Range("N125").Value = -8520
$ws.Range("H135").Value = 717.4792
$ws.Range("I135").Value = 648.907
$ws.Range("J135").Value = 1307.2
$ws.Range("K135").Value = 5840.163
$ws.Range("L135").Value = 11764.8
$ws.Range("M135").Value = -3305.163
$ws.Range("N135").Value = -16834.8
$ws.Range("H137").Value = 1206.6885
$ws.Range("I137").Value = 997.8723
$ws.Range("J137").Value = 1907.7142
$ws.Range("K137").Value = 2993.6169
$ws.Range("L137").Value = 5723.142599999999
$ws.Range("M137").Value = -443.6169
$ws.Range("N137").Value = -10823.1426

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8783.642
$ws.Range("I32").Value = 8976.286
$ws.Range("J32").Value = 7557.727
$ws.Range("K32").Value = 8976.286
$ws.Range("L32").Value = 7557.727
$ws.Range("M32").Value = -8689.286
$ws.Range("N32").Value = -8131.727
$ws.Range("H122").Value = 1202
$ws.Range("I122").Value = 1128
$ws.Range("J122").Value = 1350
$ws.Range("K122").Value = 3384
$ws.Range("L122").Value = 4050
$ws.Range("M122").Value = -934
$ws.Range("N122").Value = -8950

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 17300
$ws.Range("I22").Value = 25450
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 25450
$ws.Range("L22").Value = 1000
$ws.Range("M22").Value = -25277
$ws.Range("N22").Value = -1346
$ws.Range("H106").Value = 38514
$ws.Range("J106").Value = 38514
$ws.Range("L106").Value = 38514
$ws.Range("N106").Value = -41038
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 769
$ws.Range("I10").Value = 769
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 769
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -630
$ws.Range("N10").ClearContents()
$ws.Range("H31").Value = 2356.5518
$ws.Range("I31").Value = 1611.65
$ws.Range("J31").Value = 4011.889
$ws.Range("K31").Value = 1611.65
$ws.Range("L31").Value = 4011.889
$ws.Range("M31").Value = -1316.65
$ws.Range("N31").Value = -4601.889
$ws.Range("H34").Value = 2356.5518
$ws.Range("I34").Value = 1611.65
$ws.Range("J34").Value = 4011.889
$ws.Range("K34").Value = 1611.65
$ws.Range("L34").Value = 4011.889
$ws.Range("M34").Value = -1409.65
$ws.Range("N34").Value = -4415.889
$ws.Range("H105").Value = 8482.23
$ws.Range("I105").Value = 10757
$ws.Range("J105").Value = 899.6667
$ws.Range("K105").Value = 10757
$ws.Range("L105").Value = 899.6667
$ws.Range("M105").Value = -9010
$ws.Range("N105").Value = -4393.6667
$ws.Range("H132").Value = 797185.4
$ws.Range("I132").Value = 1230394.2
$ws.Range("J132").Value = 2969
$ws.Range("K132").Value = 3691182.6
$ws.Range("L132").Value = 8907
$ws.Range("M132").Value = -3688652.6
$ws.Range("N132").Value = -13967
$ws.Range("H134").Value = 2148.2104
$ws.Range("I134").Value = 1827.9286
$ws.Range("J134").Value = 3045
$ws.Range("K134").Value = 5483.7858
$ws.Range("L134").Value = 9135
$ws.Range("M134").Value = -2948.7858
$ws.Range("N134").Value = -14205

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 15016.5
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = 22499.75
$ws.Range("K4").Value = 150
$ws.Range("L4").Value = 67499.25
$ws.Range("M4").Value = -38
$ws.Range("N4").Value = -67723.25
$ws.Range("H70").Value = 15187.5
$ws.Range("I70").Value = 18583.334
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 55750.00199999999
$ws.Range("L70").Value = 15000
$ws.Range("M70").Value = -55435.00199999999
$ws.Range("N70").Value = -15630
$ws.Range("H73").Value = 15187.5
$ws.Range("I73").Value = 18583.334
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 55750.00199999999
$ws.Range("L73").Value = 15000
$ws.Range("M73").Value = -54658.00199999999
$ws.Range("N73").Value = -17184
$ws.Range("H122").Value = 798.13635
$ws.Range("I122").Value = 735.7143
$ws.Range("J122").Value = 907.375
$ws.Range("K122").Value = 6621.428699999999
$ws.Range("L122").Value = 8166.375
$ws.Range("M122").Value = -4171.428699999999
$ws.Range("N122").Value = -13066.375

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 21801
$ws.Range("J14").Value = 14751.25
$ws.Range("L14").Value = 14751.25
$ws.Range("N14").Value = -15095.25
$ws.Range("H75").Value = 29000
$ws.Range("J75").Value = 29000
$ws.Range("L75").Value = 29000
$ws.Range("N75").Value = -30872
$ws.Range("H78").Value = 29000
$ws.Range("J78").Value = 29000
$ws.Range("L78").Value = 87000
$ws.Range("N78").Value = -96360
$ws.Range("H80").Value = 50000
$ws.Range("J80").Value = 50000
$ws.Range("L80").Value = 50000
$ws.Range("N80").Value = -52246
$ws.Range("H83").Value = 50000
$ws.Range("J83").Value = 50000
$ws.Range("L83").Value = 150000
$ws.Range("N83").Value = -161232
$ws.Range("H86").Value = 31195
$ws.Range("J86").Value = 31195
$ws.Range("L86").Value = 31195
$ws.Range("N86").Value = -33567
$ws.Range("H87").Value = 28000
$ws.Range("J87").Value = 28000
$ws.Range("L87").Value = 28000
$ws.Range("N87").Value = -30246
$ws.Range("H88").Value = 27477.182
$ws.Range("I88").Value = 17125
$ws.Range("J88").Value = 29777.666
$ws.Range("K88").Value = 17125
$ws.Range("L88").Value = 29777.666
$ws.Range("M88").Value = -16697
$ws.Range("N88").Value = -30633.666
$ws.Range("H89").Value = 31195
$ws.Range("J89").Value = 31195
$ws.Range("L89").Value = 93585
$ws.Range("N89").Value = -105441
$ws.Range("H90").Value = 28000
$ws.Range("J90").Value = 28000
$ws.Range("L90").Value = 84000
$ws.Range("N90").Value = -95232
$ws.Range("H91").Value = 27477.182
$ws.Range("I91").Value = 17125
$ws.Range("J91").Value = 29777.666
$ws.Range("K91").Value = 17125
$ws.Range("L91").Value = 29777.666
$ws.Range("M91").Value = -15643
$ws.Range("N91").Value = -32741.666
$ws.Range("H136").Value = 1598.9656
$ws.Range("I136").Value = 1435.2174
$ws.Range("J136").Value = 2226.6667
$ws.Range("K136").Value = 4305.6522
$ws.Range("L136").Value = 6680.000100000001
$ws.Range("M136").Value = -1755.6522
$ws.Range("N136").Value = -11780.0001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 26950
$ws.Range("I29").Value = 26950
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 26950
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -26660
$ws.Range("N29").ClearContents()
$ws.Range("H31").Value = 30000
$ws.Range("I31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("H54").Value = 45666.668
$ws.Range("J54").Value = 45666.668
$ws.Range("L54").Value = 45666.668
$ws.Range("N54").Value = -46706.668
$ws.Range("H123").Value = 22353.312
$ws.Range("J123").Value = 22353.312
$ws.Range("L123").Value = 22353.312
$ws.Range("N123").Value = -32153.312
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
